$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 183 with the columns that were added to it ---
$ws.Range("D183").Value = -0.1
$ws.Range("E183").Value = -0.5
$ws.Range("M183").Value = 1.13
$ws.Range("O183").Value = 38
$ws.Range("P183").Value = 5.25
$ws.Range("Q183").Value = 1.75

# --- New row 184 ---
$ws.Range("A184").Value = "14-09-2021"
$ws.Range("B184").Value = 0.25
$ws.Range("C184").Value = 0.1
$ws.Range("D184").Value = -0.1
$ws.Range("E184").Value = -0.5
$ws.Range("F184").Value = 0.75
$ws.Range("G184").Value = 4.35
$ws.Range("H184").Value = 1.75
$ws.Range("I184").Value = 0.1
$ws.Range("J184").Value = 0.75
$ws.Range("K184").Value = 6.75
$ws.Range("L184").Value = 0.5
$ws.Range("M184").Value = 1.13
$ws.Range("N184").Value = 19
$ws.Range("O184").Value = 38
$ws.Range("P184").Value = 5.25
$ws.Range("Q184").Value = 1.75
$ws.Range("R184").Value = 4.5
$ws.Range("S184").Value = 1

# --- New row 185 ---
$ws.Range("A185").Value = "15-09-2021"
$ws.Range("B185").Value = 0.25
$ws.Range("C185").Value = 0.1
$ws.Range("D185").Value = -0.1
$ws.Range("E185").Value = -0.5
$ws.Range("F185").Value = 0.75
$ws.Range("G185").Value = 4.35
$ws.Range("H185").Value = 1.75
$ws.Range("I185").Value = 0.1
$ws.Range("J185").Value = 0.75
$ws.Range("K185").Value = 6.75
$ws.Range("L185").Value = 0.5
$ws.Range("M185").Value = 1.13
$ws.Range("N185").Value = 19
$ws.Range("O185").Value = 38
$ws.Range("P185").Value = 5.25
$ws.Range("Q185").Value = 1.75
$ws.Range("R185").Value = 4.5
$ws.Range("S185").Value = 1

# --- New row 186 (no R value for this row) ---
$ws.Range("A186").Value = "16-09-2021"
$ws.Range("B186").Value = 0.25
$ws.Range("C186").Value = 0.1
$ws.Range("D186").Value = -0.1
$ws.Range("E186").Value = -0.5
$ws.Range("F186").Value = 0.75
$ws.Range("G186").Value = 4.35
$ws.Range("H186").Value = 1.75
$ws.Range("I186").Value = 0.1
$ws.Range("J186").Value = 0.75
$ws.Range("K186").Value = 6.75
$ws.Range("L186").Value = 0.5
$ws.Range("M186").Value = 1.13
$ws.Range("N186").Value = 19
$ws.Range("O186").Value = 38
$ws.Range("P186").Value = 5.25
$ws.Range("Q186").Value = 1.75
$ws.Range("S186").Value = 1

# --- New row 187 ---
$ws.Range("A187").Value = "17-09-2021"
$ws.Range("B187").Value = 0.25
$ws.Range("C187").Value = 0.1
$ws.Range("D187").Value = -0.1
$ws.Range("E187").Value = -0.5
$ws.Range("F187").Value = 0.75
$ws.Range("G187").Value = 4.35
$ws.Range("H187").Value = 1.75
$ws.Range("I187").Value = 0.1
$ws.Range("J187").Value = 0.75
$ws.Range("K187").Value = 6.75
$ws.Range("L187").Value = 0.5
$ws.Range("M187").Value = 1.13
$ws.Range("N187").Value = 19
$ws.Range("O187").Value = 38
$ws.Range("P187").Value = 5.25
$ws.Range("Q187").Value = 1.75
$ws.Range("R187").Value = 4.5
$ws.Range("S187").Value = 1

# --- New row 188 (sparse row) ---
$ws.Range("A188").Value = "20-09-2021"
$ws.Range("B188").Value = 0.25
$ws.Range("C188").Value = 0.1
$ws.Range("I188").Value = 0.1
$ws.Range("J188").Value = 0.75
$ws.Range("K188").Value = 6.75
$ws.Range("L188").Value = 0.5
$ws.Range("N188").Value = 19
$ws.Range("R188").Value = 4.5
$ws.Range("S188").Value = 1
